# Product template update - 6/1/2026
# "cap nhat 6/1/2026 : sua edit phieu nhap xuat va sua product"
#
# Reworks the "Template" sheet:
#  - drops Mo ta / Ton kho toi thieu / Ton kho toi da / Trang thai columns
#  - adds lot / expiry / document-number / document-date tracking columns
#  - refreshes the two sample rows with the new column layout
#  - re-applies header/number formatting (currency + widths + row heights)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# ---- Header row (row 1) --------------------------------------------------
$headers = @(
    "Tên sản phẩm",
    "Mã SKU",
    "Giá bán",
    "Giá vốn",
    "Tồn kho",
    "Đơn vị",
    "Nhà cung cấp",
    "Số lô",
    "Hạn sử dụng",
    "Số chứng từ",
    "Ngày chứng từ",
    "Nhóm sản phẩm"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---- Data row 2 : Coca Cola Lon 330ml ------------------------------------
$ws.Cells.Item(2, 1).Value  = "Coca Cola Lon 330ml"
$ws.Cells.Item(2, 2).Value  = "SP000001"
$ws.Cells.Item(2, 3).Value  = 10000
$ws.Cells.Item(2, 4).Value  = 8000
$ws.Cells.Item(2, 5).Value  = 100
$ws.Cells.Item(2, 6).Value  = "Lon"
$ws.Cells.Item(2, 7).Value  = "Công ty CocaCola"
$ws.Cells.Item(2, 8).Value  = "L01"
$ws.Cells.Item(2, 9).Value  = "'2026-12-31"
$ws.Cells.Item(2, 10).Value = "NK001"
$ws.Cells.Item(2, 11).Value = "'2025-01-01"
$ws.Cells.Item(2, 12).Value = "Đồ uống"

# ---- Data row 3 : Bánh mì sandwich ----------------------------------------
$ws.Cells.Item(3, 1).Value  = "Bánh mì sandwich"
$ws.Cells.Item(3, 2).Value  = "SP000002"
$ws.Cells.Item(3, 3).Value  = 15000
$ws.Cells.Item(3, 4).Value  = 10000
$ws.Cells.Item(3, 5).Value  = 50
$ws.Cells.Item(3, 6).Value  = "Cái"
$ws.Cells.Item(3, 7).Value  = "Bánh Mỳ ABC"
$ws.Cells.Item(3, 8).Value  = "L02"
$ws.Cells.Item(3, 9).Value  = "'2026-12-31"
$ws.Cells.Item(3, 10).Value = "NK002"
$ws.Cells.Item(3, 11).Value = "'2025-01-02"
$ws.Cells.Item(3, 12).Value = "Thực Phẩm"

# ---- Column widths ---------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.8
$ws.Columns.Item(2).ColumnWidth = 15.8
$ws.Range("C1:D1").EntireColumn.ColumnWidth = 12.8
$ws.Range("E1:F1").EntireColumn.ColumnWidth = 10.8
$ws.Columns.Item(7).ColumnWidth = 20.8
$ws.Range("H1:K1").EntireColumn.ColumnWidth = 15.8

# ---- Row heights ------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30

# ---- Number formatting for the price / cost / stock columns -----------------
$ws.Range("C2:E3").NumberFormat = "$#,##0.00"
$ws.Range("C2:E3").HorizontalAlignment = -4130
$ws.Range("C2:E3").VerticalAlignment = -4130

Write-Output "done"
